# Updates cryptos list price/volume figures (and a couple of re-ranked rows)
# to match the latest GitHub Actions scrape.
#
# Note: some "Price" column values are plain numeric-looking strings
# (e.g. "7.84", "604.80") that Excel's smart-typing would otherwise coerce
# into real numbers (losing the trailing zero / exact text). Those are
# written with a leading apostrophe to force text, then the cell style is
# reset back to "Normal" so no stray number-format/quote-prefix sticks to
# the cell once the value is stored.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.240.92'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '3.550.34'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'604.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = "'144.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '3.548.55'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("D11").Value = "'7.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.97%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").Value = '4.155.77'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("D15").Value = "'30.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").Value = '3.518.99'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = '66.377.60'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = "'11.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.65%  '
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = "'14.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").Value = "'429.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D24").Value = "'79.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.88%  '
$ws.Range("D25").Value = '3.697.05'
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = "'9.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = "'7.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("B32").Value = 'RenzoRestakedETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D32").Value = '3.547.66'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'25.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("D35").Value = "'0.154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.96%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("D37").Value = "'7.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.69%  '
$ws.Range("E38").Value = '  -2.02%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = "'175.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("D41").Value = "'0.0847"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.44%  '
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = "'0.887"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.73%  '
$ws.Range("D44").Value = "'1.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").Value = "'45.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D47").Value = "'2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("D49").Value = "'25.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.22%  '
$ws.Range("E50").Value = '  -0.97%  '
$ws.Range("D51").Value = "'23.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.05%  '
